$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.345.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.10%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.689.12'
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '682.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.21%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.59'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.70%  '

# Row 7
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("E8").Value = '  -0.74%  '

# Row 9
$ws.Range("E9").Value = '  -0.96%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.36%  '

# Row 11
$ws.Range("E11").Value = '  -0.38%  '

# Row 12
$ws.Range("E12").Value = '  -2.74%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.311.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.16%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.51'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.53%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.691.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.04%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.343.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.14%  '

# Row 17
$ws.Range("E17").Value = '  +1.94%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '16.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.48%  '

# Row 19
$ws.Range("E19").Value = '  -1.50%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.62'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.21%  '

# Row 21
$ws.Range("E21").Value = '  -0.58%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.656'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.01%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.835.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.14%  '

# Row 25
$ws.Range("E25").Value = '  -0.07%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.79%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.71%  '

# Row 28
$ws.Range("E28").Value = '  -3.72%  '

# Row 29
$ws.Range("E29").Value = '  -0.69%  '

# Row 30
$ws.Range("E30").Value = '  -4.30%  '

# Row 31
$ws.Range("E31").Value = '  -2.90%  '

# Row 32
$ws.Range("E32").Value = '  -3.40%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.95'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.01%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.677.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.29%  '

# Row 36
$ws.Range("E36").Value = '  -6.81%  '

# Row 37
$ws.Range("E37").Value = '  -2.22%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.29'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.40%  '

# Row 39
$ws.Range("E39").Value = '  -0.01%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.25'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.62%  '

# Row 41
$ws.Range("E41").Value = '  -0.11%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0907'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.27%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '169.95'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.73%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.944'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.08%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '47.64'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.93%  '

# Row 46
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.03%  '

# Row 47
$ws.Range("E47").Value = '  -2.77%  '

# Row 48
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '28.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.30%  '

# Row 49
$ws.Range("B49").Value = 'FLOKI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000278'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.88%  '

# Row 50
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.56%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.10%  '
